$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5: J5 corrected value
$ws.Range("J5").Value = 535.02

# Row 8 (n=5)
$ws.Range("B8").Value = 240.49
$ws.Range("C8").Value = 231.17
$ws.Range("D8").Value = 257.92
$ws.Range("E8").Value = 253.47
$ws.Range("F8").Value = 1168.83
$ws.Range("G8").Value = 1062.18
$ws.Range("H8").Value = 1144.1199999999999
$ws.Range("I8").Value = 1202.8800000000001
$ws.Range("J8").Value = 2987.73
$ws.Range("K8").Value = 2875.32
$ws.Range("L8").Value = 3083.99
$ws.Range("M8").Value = 2996.66
$ws.Range("N8").Value = 4910.5200000000004
$ws.Range("O8").Value = 4741.32
$ws.Range("P8").Value = 5048.5
$ws.Range("Q8").Value = 4774.8500000000004

# Row 9 (n=6)
$ws.Range("B9").Value = 268.05
$ws.Range("C9").Value = 265.73
$ws.Range("D9").Value = 280
$ws.Range("E9").Value = 280
$ws.Range("F9").Value = 1342.79
$ws.Range("G9").Value = 1310.82
$ws.Range("H9").Value = 1353.5
$ws.Range("I9").Value = 1348.58
$ws.Range("J9").Value = 3415.66
$ws.Range("K9").Value = 3227.09
$ws.Range("L9").Value = 3466.5
$ws.Range("M9").Value = 3359.1
$ws.Range("N9").Value = 5584.42
$ws.Range("O9").Value = 5370.85
$ws.Range("P9").Value = 5683.93
$ws.Range("Q9").Value = 5419.2
$ws.Range("R9").Value = 6917.29
$ws.Range("S9").Value = 6728.85
$ws.Range("T9").Value = 6984.15
$ws.Range("U9").Value = 6677.03

# Row 10 (n=7)
$ws.Range("B10").Value = 288
$ws.Range("C10").Value = 288
$ws.Range("D10").Value = 288
$ws.Range("E10").Value = 288
$ws.Range("F10").Value = 1422.85
$ws.Range("G10").Value = 1420.61
$ws.Range("H10").Value = 1419.09
$ws.Range("I10").Value = 1418.03
$ws.Range("J10").Value = 3579.98
$ws.Range("K10").Value = 3508.94
$ws.Range("L10").Value = 3612.09
$ws.Range("M10").Value = 3538.88
$ws.Range("N10").Value = 5958.78
$ws.Range("O10").Value = 5840.98
$ws.Range("P10").Value = 5974.08
$ws.Range("Q10").Value = 5786.8
$ws.Range("R10").Value = 7407.25
$ws.Range("S10").Value = 7298.75
$ws.Range("T10").Value = 7388.34
$ws.Range("U10").Value = 7073.54

# Row 11 (n=8)
$ws.Range("B11").Value = 288
$ws.Range("C11").Value = 288
$ws.Range("D11").Value = 288
$ws.Range("E11").Value = 282.25
$ws.Range("F11").Value = 1457.73
$ws.Range("G11").Value = 1464.03
$ws.Range("H11").Value = 1449.91
$ws.Range("I11").Value = 1449.82
$ws.Range("J11").Value = 3665.77
$ws.Range("K11").Value = 3630.66
$ws.Range("L11").Value = 3666.91
$ws.Range("M11").Value = 3628.95
$ws.Range("N11").Value = 6075.37
$ws.Range("O11").Value = 6009.09
$ws.Range("P11").Value = 6087.34
$ws.Range("Q11").Value = 5950.3
$ws.Range("R11").Value = 7589.69
$ws.Range("S11").Value = 7532.33
$ws.Range("T11").Value = 7574.47
$ws.Range("U11").Value = 7401.15
